$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edit order matters for shared-string table ordering; replicate the
# order the original author most likely used when retyping headers
# (new data cell first, then new/renamed headers from right to left).
$ws.Range("C2").Value = "Ethan"
$ws.Range("L2").Value = 90

$ws.Range("K1").Value = "CARBON_FOOT_PRINT"
$ws.Range("I1").Value = "CAR_TYPE"
$ws.Range("H1").Value = "CAR_FUEL_TYPE"
$ws.Range("G1").Value = "GAS_ACCOUNT_PAYMENT"
$ws.Range("F1").Value = "WATER ACCOUNT_PAYMENT"
$ws.Range("E1").Value = "ELECTRICAL_ACCOUNT_PAYMENT"
$ws.Range("B1").Value = "ID_NUMBER"
$ws.Range("A1").Value = "FULL_NAME"

# Unchanged headers (kept for completeness / no-ops for the string table)
$ws.Range("C1").Value = "USERNAME"
$ws.Range("D1").Value = "PASSWORD"
$ws.Range("J1").Value = "DISTANCE"

# Selection / view state
$ws.Range("C10").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
